$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day block: header row (styled like the previous header row) + 35 state/UT data rows.
$headerRow = 1045
$startRow = 1046
$date = "16-10-2020"

$headers = @("States/UT", "Active Cases", "Active Cases Since Yesterday", "Recovered Cases", "Recovered Cases Since Yesterday", "Deceased Cases", "Deceased Cases Since Yesterday", "Date")
for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item($headerRow, $c).Value = $headers[$c - 1]
}

# Match formatting of the previous block header row (row 1009).
$ws.Range("A1009:H1009").Copy()
$ws.Range("A1045:H1045").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @("Andaman and Nicobar Islands", 190, -5, 3817, 21, 55, 0),
    @("Andhra Pradesh", 40047, -1622, 725099, 5622, 6357, 38),
    @("Arunachal Pradesh", 3052, 7, 9889, 195, 30, 1),
    @("Assam", 28804, -503, 169335, 1263, 843, 9),
    @("Bihar", 11038, 282, 189186, 806, 972, 5),
    @("Chandigarh", 1044, -41, 12232, 113, 201, 2),
    @("Chhattisgarh", 28187, 378, 123943, 2395, 1385, 46),
    @("Dadra and Nagar Haveli and Daman and Diu", 71, -8, 3099, 9, 2, 0),
    @("Delhi", 22605, 702, 292502, 2755, 5924, 26),
    @("Goa", 4084, -104, 35161, 430, 525, 6),
    @("Gujarat", 14782, -155, 137733, 1329, 3606, 11),
    @("Haryana", 10364, 177, 134719, 1013, 1623, 9),
    @("Himachal Pradesh", 2654, 134, 15389, 156, 260, 5),
    @("Jammu and Kashmir", 9058, -681, 75641, 1323, 1358, 6),
    @("Jharkhand", 6892, -299, 87240, 873, 820, 9),
    @("Karnataka", 113557, -449, 620008, 8841, 10283, 85),
    @("Kerala", 94609, 684, 222231, 7082, 1089, 23),
    @("Ladakh", 1018, 39, 4310, 49, 65, 1),
    @("Madhya Pradesh", 14157, -275, 139717, 1559, 2710, 24),
    @("Maharashtra", 192936, -3825, 1330483, 13714, 41196, 337),
    @("Manipur", 3193, 126, 11081, 166, 104, 1),
    @("Meghalaya", 2445, 106, 5646, 64, 73, 3),
    @("Mizoram", 108, -4, 2121, 13, 0, 0),
    @("Nagaland", 1453, -25, 6017, 101, 22, 0),
    @("Odisha", 22387, -329, 238535, 2772, 1089, 27),
    @("Puducherry", 4551, 26, 27365, 213, 570, 2),
    @("Punjab", 7090, -670, 115186, 1111, 3954, 29),
    @("Rajasthan", 21587, -124, 143984, 2149, 1708, 14),
    @("Sikkim", 312, -13, 3129, 54, 59, 0),
    @("Tamil Nadu", 41872, -694, 622458, 5055, 10472, 49),
    @("Telengana", 23315, 112, 194653, 1435, 1256, 7),
    @("Tripura", 3105, -213, 25765, 371, 323, 4),
    @("Uttarakhand", 5682, -463, 49997, 868, 814, 18),
    @("Uttar Pradesh", 36295, -603, 404545, 3239, 6543, 36),
    @("West Bengal", 31984, 479, 271563, 3179, 5870, 62),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $date
}
